# feat: PostgreSQL session 02 despues del taller
#
# 1) Slide 14 ("Eliminacion de tablas"): the TRUNCATE TABLE example box
#    gains the same "[CASCADE | RESTRICT]" clause already shown on the
#    DROP TABLE example box above it, and is repositioned/resized to fit.
# 2) Slide 7: numeric(s,p) -> numeric(p,s) typo fix.

$p = $ppt.ActivePresentation

# Helper: convert an EMU integer to points, nudged by a hair so the
# engine's point->EMU round-trip lands back on the exact EMU value
# instead of one EMU short (floating point truncation).
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + 0.00001
}

# ---------------------------------------------------------------------
# 1) Slide 14, shape "CuadroTexto 6" (the TRUNCATE TABLE ...; box)
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$shape = $s14.Shapes.Item(6)

# Reposition / resize the box to fit the longer text.
$shape.Left   = EmuToPt 2567729
$shape.Top    = EmuToPt 4139223
$shape.Width  = EmuToPt 7056539
$shape.Height = EmuToPt 369332

$tr = $shape.TextFrame.TextRange

# Original text: "TRUNCATE TABLE nombre_tabla;" (28 chars). Replace the
# trailing ";" with " [CASCADE | RESTRICT];".
$semi = $tr.Characters(28, 1)
$semi.Text = " [CASCADE | RESTRICT];"

# Now color/segment the newly inserted text to match the DROP TABLE
# example's styling: " [" stays default (inherited from the old ";"
# run: Courier New / E8E6E3), "CASCADE" is highlighted orange, " | "
# + "RESTRICT" stay E8E6E3, "];" stays default too.
$cascade = $tr.Characters(30, 7)
$cascade.Font.Color.RGB = 8239868   # FCBA7D (orange) in BGR

$restrict = $tr.Characters(37, 11)
$restrict.Font.Color.RGB = 14935784 # E8E6E3 in BGR

# ---------------------------------------------------------------------
# 2) Slide 7, shape "CuadroTexto 10": numeric(s,p) -> numeric(p,s)
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shape2 = $s7.Shapes.Item(5)
$tr2 = $shape2.TextFrame.TextRange
$frag = $tr2.Characters(37, 3)
$frag.Text = "p,s"
